# Apply the "rough runnable code for backprojection_MCNP_Gammas" edit:
#  - Add a new worksheet "工作表1" after Sheet1, populated with a lookup/interpolation
#    table (A1:P10) of values used for backprojection weighting.
#  - Restore Sheet1 as the active sheet, with selection moved to G2.
#  - Select A1:P1 on the new sheet (matches the header-row highlight left behind
#    by the author when they finished pasting the table).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new sheet right after Sheet1 --------------------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "工作表1"

$sheet2Data = @(
    @(-2.5,0.58713351915171197,0.48291820425854498,0.43939861829540899,0.39721456856431298,0.36983940712879998,0.342752590385347,0.32228974774039199,0.30425488535119199,0.28844913104541398,0.27234731623875102,0.26286140635437599,0.25254163390884499,0.24503555307861599,0.23917541804407799,0.23142241571064301),
    @(-2,0.51089037336028198,0.42020810195367903,0.38233982021545598,0.34563364655310802,0.32181332972501497,0.29824391413586898,0.28043830607929299,0.26474538908686501,0.25099211584471298,0.236981227503268,0.22872712535374001,0.219747443176323,0.213216076267875,0.208116918277266,0.20137069424543999),
    @(-1.5,0.45443314139090901,0.37377194358300903,0.34008839202929603,0.30743852686133799,0.28625053434952102,0.26528571660109501,0.24944776226599,0.23548903072922001,0.223255597707708,0.21079301799456501,0.203451056265084,0.195463696562579,0.189654094862197,0.18511843221332,0.17911771671902499),
    @(-1,0.41724313482505099,0.34318310714066902,0.31225615802050599,0.28227830021595002,0.26282429563076198,0.24357520158006099,0.229033397485912,0.216217023911668,0.20498475346601999,0.193542089289716,0.18680098075525101,0.17946729247924201,0.17413313833253299,0.16996866631287799,0.16445903878797999),
    @(-0.5,0.41435816817981502,0.34081021772758202,0.31009710847488198,0.28032652818466602,0.26100703546956,0.241891036465014,0.227449779548929,0.21472202291571699,0.20356741636164899,0.19220387081361301,0.18550937268840401,0.178226392127627,0.17292912020979001,0.16879344282297301,0.16332191081196301),
    @(0,0.44074518254712702,0.36251357680750301,0.32984460589367998,0.29817818574759197,0.27762839574142101,0.25729505826230298,0.24193415818966801,0.228395876935645,0.216530926563096,0.20444373161541801,0.19732291676285299,0.189576143938638,0.183941532972273,0.17954248880037799,0.17372252056951401),
    @(0.5,0.49712674385147798,0.408887495942163,0.37203940371456201,0.33632211181886601,0.31314352564905101,0.29020908131953999,0.27288316480786001,0.25761302245879703,0.244230268935373,0.230596840585888,0.222565107871836,0.21382734260074801,0.20747193382148299,0.20251014957101299,0.195945672021406),
    @(1,0.56846877102151405,0.46756642079539701,0.42543030568554502,0.38458727062600201,0.35808235501456498,0.33185661772863601,0.31204428095932601,0.29458273988985401,0.27927944441754798,0.26368950009349701,0.25450514362585303,0.24451343276627099,0.23724596735071601,0.23157212374880901,0.22406558636940299),
    @(1.5,0.66067490545242202,0.54340610531099098,0.49443547528621501,0.44696766403272298,0.41616362780722599,0.38568405287707502,0.36265813766569999,0.34236431928338801,0.32457881583153098,0.30646016883216598,0.29578610167101399,0.28417372652576101,0.27572747183058699,0.26913332580827098,0.26040922146652801),
    @(2,0.72383687498744997,0.59535692043283495,0.54170459080631905,0.489698752720653,0.45594978309209599,0.42255629397012201,0.39732905077868602,0.37509509886347497,0.35560926228571799,0.33575843290687302,0.32406390152150999,0.31134135785146499,0.30208762268845502,0.29486306184830202,0.28530491400330699)
)

for ($r = 1; $r -le $sheet2Data.Count; $r++) {
    $rowVals = $sheet2Data[$r - 1]
    for ($c = 1; $c -le $rowVals.Count; $c++) {
        $ws2.Cells.Item($r, $c).Value = $rowVals[$c - 1]
    }
}

# Leave the new sheet's header row selected, as in the saved file.
[void]$ws2.Range("A1:P1").Select()

# --- Restore Sheet1 as the active sheet, selection on G2 -----------------
$ws1.Activate()
[void]$ws1.Range("G2").Select()
